# Generate Report for Handback
#
# This script replays (via Excel COM automation) the "handback" report
# generation that happens once a localized target file has been produced
# and is back "in sync" with the en-US source:
#   - the Status column moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" (shared across Overview + every
#     language sheet, since they all pointed at the same cell text)
#   - the per-language sheets get their "Latest Target File" / "Latest
#     Handback File" / "Latest Handback DateTime" columns filled in
#   - the "Latest Target File" cell becomes a hyperlink to the source .md
#   - a handful of columns get widened so the new, longer content fits

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Status text: update every cell that shows it so the old shared
#    string is fully replaced (Overview!E2/F2, and each language
#    sheet's Status cell C2).
# ---------------------------------------------------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn handback details
# ---------------------------------------------------------------------
$wsZhCn.Range("J2").Value = "0f51c9b3-f0a2-4be6-b663-558a8ba7405b.b53ab6119df9594e81535c4c4b1de34c79abd09a.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-30 19:13:48"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1811585a046efea986da45231734b31b65e72af5/e2e/0f51c9b3-f0a2-4be6-b663-558a8ba7405b.md",
    "",
    "",
    "0f51c9b3-f0a2-4be6-b663-558a8ba7405b.md"
) | Out-Null

# ---------------------------------------------------------------------
# 3. de-de handback details
# ---------------------------------------------------------------------
$wsDeDe.Range("J2").Value = "0f51c9b3-f0a2-4be6-b663-558a8ba7405b.b53ab6119df9594e81535c4c4b1de34c79abd09a.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-30 19:13:55"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1811585a046efea986da45231734b31b65e72af5/e2e/0f51c9b3-f0a2-4be6-b663-558a8ba7405b.md",
    "",
    "",
    "0f51c9b3-f0a2-4be6-b663-558a8ba7405b.md"
) | Out-Null

# ---------------------------------------------------------------------
# 4. Widen the columns that now hold the longer handback text/links.
#    ColumnWidth is snapped by Excel to a whole-pixel grid, so we feed
#    it the character-width value that lands exactly on the target
#    stored width (~29.98 and 40 "characters" respectively).
# ---------------------------------------------------------------------
$wideWidth  = 29.166666666666668   # -> stored width 30 (closest pixel match to 29.9777047293527)
$fullWidth  = 39.166666666666664   # -> stored width 40

$wsOverview.Columns.Item(5).ColumnWidth = $wideWidth   # E: zh-cn
$wsOverview.Columns.Item(6).ColumnWidth = $wideWidth   # F: de-de

$wsZhCn.Columns.Item(3).ColumnWidth = $wideWidth    # C: Status
$wsZhCn.Columns.Item(9).ColumnWidth = $fullWidth    # I: Latest Target File
$wsZhCn.Columns.Item(10).ColumnWidth = $fullWidth   # J: Latest Handback File

$wsDeDe.Columns.Item(3).ColumnWidth = $wideWidth    # C: Status
$wsDeDe.Columns.Item(9).ColumnWidth = $fullWidth    # I: Latest Target File
$wsDeDe.Columns.Item(10).ColumnWidth = $fullWidth   # J: Latest Handback File
